{"js": "// Update the date and the 25 two-digit multiplication problems to the new\n// values. Every \"before\" string is unique within the document, so a\n// plain case-sensitive search/replace for each pair is unambiguous.\nconst replacements = [\n  [\"2024-04-05 Friday\", \"2024-04-06 Saturday\"],\n  [\"91\u00d761=5551\", \"37\u00d746=1702\"],\n  [\"27\u00d772=1944\", \"26\u00d740=1040\"],\n  [\"35\u00d742=1470\", \"86\u00d716=1376\"],\n  [\"92\u00d788=8096\", \"55\u00d767=3685\"],\n  [\"21\u00d781=1701\", \"83\u00d713=1079\"],\n  [\"34\u00d796=3264\", \"74\u00d734=2516\"],\n  [\"68\u00d747=3196\", \"13\u00d726=338\"],\n  [\"37\u00d766=2442\", \"49\u00d798=4802\"],\n  [\"95\u00d720=1900\", \"47\u00d795=4465\"],\n  [\"90\u00d773=6570\", \"24\u00d794=2256\"],\n  [\"78\u00d758=4524\", \"73\u00d747=3431\"],\n  [\"69\u00d784=5796\", \"50\u00d718=900\"],\n  [\"88\u00d727=2376\", \"82\u00d761=5002\"],\n  [\"38\u00d781=3078\", \"39\u00d799=3861\"],\n  [\"46\u00d721=966\", \"73\u00d727=1971\"],\n  [\"84\u00d777=6468\", \"95\u00d790=8550\"],\n  [\"30\u00d784=2520\", \"92\u00d796=8832\"],\n  [\"65\u00d766=4290\", \"50\u00d734=1700\"],\n  [\"13\u00d787=1131\", \"77\u00d748=3696\"],\n  [\"70\u00d792=6440\", \"56\u00d720=1120\"],\n  [\"89\u00d713=1157\", \"67\u00d744=2948\"],\n  [\"17\u00d723=391\", \"29\u00d789=2581\"],\n  [\"41\u00d729=1189\", \"62\u00d724=1488\"],\n  [\"18\u00d789=1602\", \"14\u00d781=1134\"],\n  [\"43\u00d770=3010\", \"13\u00d796=1248\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date and the 25 two-digit multiplication problems to the new\n# values. Every \"before\" string is unique within the document, so a plain\n# case-sensitive Find/Replace for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-04-05 Friday\", \"2024-04-06 Saturday\"),\n    @(\"91\u00d761=5551\", \"37\u00d746=1702\"),\n    @(\"27\u00d772=1944\", \"26\u00d740=1040\"),\n    @(\"35\u00d742=1470\", \"86\u00d716=1376\"),\n    @(\"92\u00d788=8096\", \"55\u00d767=3685\"),\n    @(\"21\u00d781=1701\", \"83\u00d713=1079\"),\n    @(\"34\u00d796=3264\", \"74\u00d734=2516\"),\n    @(\"68\u00d747=3196\", \"13\u00d726=338\"),\n    @(\"37\u00d766=2442\", \"49\u00d798=4802\"),\n    @(\"95\u00d720=1900\", \"47\u00d795=4465\"),\n    @(\"90\u00d773=6570\", \"24\u00d794=2256\"),\n    @(\"78\u00d758=4524\", \"73\u00d747=3431\"),\n    @(\"69\u00d784=5796\", \"50\u00d718=900\"),\n    @(\"88\u00d727=2376\", \"82\u00d761=5002\"),\n    @(\"38\u00d781=3078\", \"39\u00d799=3861\"),\n    @(\"46\u00d721=966\", \"73\u00d727=1971\"),\n    @(\"84\u00d777=6468\", \"95\u00d790=8550\"),\n    @(\"30\u00d784=2520\", \"92\u00d796=8832\"),\n    @(\"65\u00d766=4290\", \"50\u00d734=1700\"),\n    @(\"13\u00d787=1131\", \"77\u00d748=3696\"),\n    @(\"70\u00d792=6440\", \"56\u00d720=1120\"),\n    @(\"89\u00d713=1157\", \"67\u00d744=2948\"),\n    @(\"17\u00d723=391\", \"29\u00d789=2581\"),\n    @(\"41\u00d729=1189\", \"62\u00d724=1488\"),\n    @(\"18\u00d789=1602\", \"14\u00d781=1134\"),\n    @(\"43\u00d770=3010\", \"13\u00d796=1248\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
